$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Sheet1" (sheet4.xml) -- rebuilt as the "HR lo/hi 95" summary
# for relevel(lapatrasno, ref = "2") instead of ref = "3"
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Drop the old blank helper block D2:G3 (also shrinks used range back
# down from G to F)
$ws4.Range("D2:G3").Clear()

# New headers (reuse existing D1/F1 cells, touch E1 so it still shows
# up as a bare cell in the xml like it used to after the old D1 text)
$ws4.Range("D1").Value = "HR lo 95"
$ws4.Range("F1").Value = "HR hi 95"
$ws4.Range("E1").ClearFormats()

# Row 2 + 3: new relevel(...) labels/values; switch from the old
# Lucida-console style (s=3) to the plain 0.000-formatted style (s=1)
$ws4.Range("A2").ClearFormats()
$ws4.Range("A2").NumberFormat = "0.000"
$ws4.Range("A2").Value = "relevel(lapatrasno, ref = ""2"")1"

$ws4.Range("A3").ClearFormats()
$ws4.Range("A3").NumberFormat = "0.000"
$ws4.Range("A3").Value = "relevel(lapatrasno, ref = ""2"")3"

$ws4.Range("C2").Value = "(0.763,1.338)"
$ws4.Range("C3").Value = "(1.361,3.190)"
$ws4.Range("B2").Value = 1.010705
$ws4.Range("B3").Value = 2.0837330000000001

# Column widths for D:F + portrait page setup
$ws4.Columns.Item(4).ColumnWidth = 9.140625
$ws4.Columns.Item(5).ColumnWidth = 9.140625
$ws4.Columns.Item(6).ColumnWidth = 9.140625
$ws4.PageSetup.Orientation = 1

$ws4.Range("B2:C2").Select()

# ------------------------------------------------------------------
# Sheet "lapat only treat" (sheet3.xml) -- add the "Lapat vs Trastuz"
# row underneath the existing lapat/trastuz-vs-none rows
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A5").Value = "Lapat vs Trastuz"
$ws3.Range("B5").Value = 0.95799999999999996
$ws3.Range("J5").Value = "(0.673,1.28)"
$ws3.Range("I5").Value = 0.92700000000000005
$ws3.Range("C5").Value = "(0.693,1.324)"
$ws3.Range("E5").Value = 1.0089999999999999
$ws3.Range("F5").Value = "(0.680,1.496)"

$ws3.Range("N5").Value = "(0.763,1.338)"
$ws3.Range("N5").NumberFormat = "0.000"
$ws3.Range("L5").Value = "Lapat vs Trastuz"
$ws3.Range("M5").Value = 1.010705
$ws3.Range("M5").NumberFormat = "0.000"

$ws3.Columns.Item(12).ColumnWidth = 11.7109375

$ws3.Range("I16").Select()
